$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 11-13 (data shrinks from 13 data-bearing rows down to 9)
$ws.Range("A11:H13").EntireRow.Delete()

# Row 2
$ws.Cells.Item(2,1).Value = 2
$ws.Cells.Item(2,2).Value = "life-dev/main"
$ws.Cells.Item(2,3).Value = ""
$ws.Cells.Item(2,4).Value = "impression"
$ws.Cells.Item(2,5).Value = ""
$ws.Cells.Item(2,6).Value = "channel, page_url, banner_text, banner_position, os_name, impression_type"
$ws.Cells.Item(2,7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, ______포인트 12,500원 놓치고 있어요!___전국 날씨특파원, 오늘 날씨는?___여름맞이 체력 증진! 오운완 챌린지___, 라이프 메인 상단 카드 배너, iOS, 메인 상단 카드 배너"
$ws.Cells.Item(2,8).Value = 6

# Row 3
$ws.Cells.Item(3,1).Value = 3
$ws.Cells.Item(3,2).Value = "life-dev/main"
$ws.Cells.Item(3,3).Value = ""
$ws.Cells.Item(3,4).Value = "pageview"
$ws.Cells.Item(3,5).Value = ""
$ws.Cells.Item(3,6).Value = "channel, page_url, os_name"
$ws.Cells.Item(3,7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, iOS"
$ws.Cells.Item(3,8).Value = 3

# Row 4
$ws.Cells.Item(4,1).Value = 4
$ws.Cells.Item(4,2).Value = "life-dev/main"
$ws.Cells.Item(4,3).Value = ""
$ws.Cells.Item(4,4).Value = "swipe"
$ws.Cells.Item(4,5).Value = ""
$ws.Cells.Item(4,6).Value = "channel, page_url, swipe_area, swipe_direct, os_name"
$ws.Cells.Item(4,7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, 메인 상단 카드 배너, left, iOS"
$ws.Cells.Item(4,8).Value = 5

# Row 5
$ws.Cells.Item(5,1).Value = 8
$ws.Cells.Item(5,2).Value = "ecommerce-dev/product/detail/800"
$ws.Cells.Item(5,3).Value = ""
$ws.Cells.Item(5,4).Value = "click"
$ws.Cells.Item(5,5).Value = ""
$ws.Cells.Item(5,6).Value = "channel, page_url, tab_name, prd_code, prd_name, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, os_name"
$ws.Cells.Item(5,7).Value = "Rround, https://ecommerce-dev.hectoinnovation.co.kr/product/detail/800, 상품상세`n, 800, 여성용 스킨핏 50수 투톤 모달 팬티 5P SET, 20,000원, 20,000원, 10%, 0, 0, #여성팬티___#50수팬티___#숙녀팬티___#여자팬티___#팬티세트___#모달팬티___#투톤팬티___#팬티, iOS"
$ws.Cells.Item(5,8).Value = 12

# Row 6
$ws.Cells.Item(6,1).Value = 9
$ws.Cells.Item(6,2).Value = "life-dev/main"
$ws.Cells.Item(6,3).Value = "상품 찜하기"
$ws.Cells.Item(6,4).Value = "click"
$ws.Cells.Item(6,5).Value = ""
$ws.Cells.Item(6,6).Value = "channel, page_url, prd_code, prd_name, prd_brand, prd_price_final, prd_is_ad, os_name"
$ws.Cells.Item(6,7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, 800, 여성용 스킨핏 50수 투톤 모달 팬티 5P SET, 마이그스토어, 20,000원, F, iOS"
$ws.Cells.Item(6,8).Value = 8

# Row 7
$ws.Cells.Item(7,1).Value = 10
$ws.Cells.Item(7,2).Value = "life-dev/main"
$ws.Cells.Item(7,3).Value = ""
$ws.Cells.Item(7,4).Value = "click"
$ws.Cells.Item(7,5).Value = "상품 더보기"
$ws.Cells.Item(7,6).Value = "channel, page_url, click_text, module_id, module_order, module_name, os_name"
$ws.Cells.Item(7,7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, 상품 더보기, C-3, 33, commerce-category-ranking, iOS"
$ws.Cells.Item(7,8).Value = 7

# Row 8
$ws.Cells.Item(8,1).Value = 11
$ws.Cells.Item(8,2).Value = "ecommerce-dev/category/detail/543"
$ws.Cells.Item(8,3).Value = ""
$ws.Cells.Item(8,4).Value = "pageview"
$ws.Cells.Item(8,5).Value = ""
$ws.Cells.Item(8,6).Value = "channel, page_url, ctgr_id, os_name"
$ws.Cells.Item(8,7).Value = "Rround, https://ecommerce-dev.hectoinnovation.co.kr/category/detail/543, 543, iOS"
$ws.Cells.Item(8,8).Value = 4

# Row 9
$ws.Cells.Item(9,1).Value = 12
$ws.Cells.Item(9,2).Value = "life-dev/main"
$ws.Cells.Item(9,3).Value = "상품"
$ws.Cells.Item(9,4).Value = "click"
$ws.Cells.Item(9,5).Value = "드시모네 베이비스텝2 100억 생유산균 2박스"
$ws.Cells.Item(9,6).Value = "channel, page_url, click_text, module_id, module_order, prd_order, prd_code, prd_name, prd_brand, prd_price_final, prd_is_ad, el_order, module_name, os_name"
$ws.Cells.Item(9,7).Value = "Rround, https://life-dev.hectoinnovation.co.kr/main, 드시모네 베이비스텝2 100억 생유산균 2박스, C-3, 33, 2, 1030, 드시모네 베이비스텝2 100억 생유산균 2박스, 마이그스토어, 96,000원, F, 2, commerce-category-ranking, iOS"
$ws.Cells.Item(9,8).Value = 14

# Row 10
$ws.Cells.Item(10,1).Value = 13
$ws.Cells.Item(10,2).Value = "ecommerce-dev/product/detail/1030"
$ws.Cells.Item(10,3).Value = ""
$ws.Cells.Item(10,4).Value = "pageview"
$ws.Cells.Item(10,5).Value = ""
$ws.Cells.Item(10,6).Value = "channel, page_url, prd_code, prd_price_origin, prd_price_final, prd_disc_rate, prd_review_cnt, prd_review_score, prd_tag, os_name"
$ws.Cells.Item(10,7).Value = "Rround, https://ecommerce-dev.hectoinnovation.co.kr/product/detail/1030, 1030, 96,000원, 86,400원, 10%, 0, 0, #프로바이오틱스___#식품___#영양제___#드시모네___#베이비스텝___#박스___#생유산균___#건강식품, iOS"
$ws.Cells.Item(10,8).Value = 10
